# Fix duplicate word "discount discount." -> "discount." on the Module 2
# "Question #1" slide (Content Placeholder 2, bullet "Yes: the number of
# products ordered was higher with a discount discount.").
#
# Locate the run by its known text rather than hard-coding slide/shape
# indices so the edit is resilient to minor structural differences.

$p = $ppt.ActivePresentation

$needle = "a discount discount."
$replacement = "a discount."

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)

        if ($shp.HasTextFrame -eq $false) { continue }
        if ($shp.TextFrame.HasText -eq $false) { continue }

        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        $idx = $fullText.IndexOf($needle)

        if ($idx -ge 0) {
            $startPos = $idx + 1
            $target = $tr.Characters($startPos, $needle.Length)
            $target.Text = $replacement
        }
    }
}
